$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("AddEmployee")

# New "Employee ID" column header
$ws.Range("D1").Value = "Employee ID"

# Row 2 - John / J / Doe (unchanged), add Employee ID
$ws.Range("D2").Value = 222222

# Row 3 - was Mike / M / Smith -> Alijon / Asel / John
$ws.Range("A3").Value = "Alijon"
$ws.Range("B3").Value = "Asel"
$ws.Range("C3").Value = "John"
$ws.Range("D3").Value = 333333

# Row 4 - was Ryan / R / Carter -> Donald / R / Trump
$ws.Range("A4").Value = "Donald"
$ws.Range("B4").Value = "R"
$ws.Range("C4").Value = "Trump"
$ws.Range("D4").Value = 555555

# Row 5 - was James / J / John -> Joe / J / Bidan
$ws.Range("A5").Value = "Joe"
$ws.Range("B5").Value = "J"
$ws.Range("C5").Value = "Bidan"
$ws.Range("D5").Value = 777777

# Row 6 - Donald / D / Duck (unchanged), add Employee ID
$ws.Range("D6").Value = 8888888

# New column was auto-fit to its (wider) content, like columns A-C already are
$ws.Columns.Item(4).ColumnWidth = 18.166666666666668

# Match the final selection recorded in the saved workbook
[void]$ws.Range("D6").Select()
